# DOMA-3100 add formatter convert to number for some colomns
#
# Adds a ":formatN()" modifier to the ticket-count placeholders (processing,
# completed, canceled, deferred, closed, new_or_reopened) for both the
# {d.tickets[i].*} and {d.tickets[i+1].*} rows, and applies a numeric
# ("0") number format to those same cells so the generated report renders
# the counts as numbers instead of plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the six "count" columns (address stays in column A)
$countCols = @("B", "C", "D", "E", "F", "G")

# Row 2 -> {d.tickets[i].*}, row 3 -> {d.tickets[i+1].*}
$exprs2 = @(
    "{d.tickets[i].processing:formatN()}",
    "{d.tickets[i].completed:formatN()}",
    "{d.tickets[i].canceled:formatN()}",
    "{d.tickets[i].deferred:formatN()}",
    "{d.tickets[i].closed:formatN()}",
    "{d.tickets[i].new_or_reopened:formatN()}"
)

$exprs3 = @(
    "{d.tickets[i+1].processing:formatN()}",
    "{d.tickets[i+1].completed:formatN()}",
    "{d.tickets[i+1].canceled:formatN()}",
    "{d.tickets[i+1].deferred:formatN()}",
    "{d.tickets[i+1].closed:formatN()}",
    "{d.tickets[i+1].new_or_reopened:formatN()}"
)

for ($i = 0; $i -lt $countCols.Length; $i++) {
    $col = $countCols[$i]
    $ws.Range($col + "2").Value = $exprs2[$i]
    $ws.Range($col + "3").Value = $exprs3[$i]
}

# Apply a plain integer number format to the updated count cells so Excel
# assigns a distinct numFmtId (builtin "0") rather than the text "@" format
# used by the header/address cells.
$ws.Range("B2:G2").NumberFormat = "0"
$ws.Range("B3:G3").NumberFormat = "0"
